# Apply "new reports in excel sheet" changes to TcRAReportList.xlsx
#
# Sheet1 = ProjectReports        (A:ReportName, B:ProjectValue, C:ExpectedRows)
# Sheet2 = MultiProjectReports   (A:ReportName, B:FilterKey,   C:FilterValue, D:ExpectedRows)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Helper: write a value into a cell that keeps it a true NUMBER even though
# the cell's style carries a Text ("@") number format (numFmtId 49). Writing
# directly to such a cell would otherwise coerce the value into a text
# shared-string, same as typing a number into a Text-formatted Excel cell.
function Set-NumericValue($range, $number) {
    $range.NumberFormat = "General"
    $range.Value = $number
    $range.NumberFormat = "@"
}

# ---------------------------------------------------------------------
# Sheet1 (ProjectReports): replace the PROJECT_13062018 rows with the new
# 312498* rows (one additional row appended).
# ---------------------------------------------------------------------

$ws1.Range("A2").Value = "Comment Status Internal (Project)"
$ws1.Range("B2").Value = "312498*"
Set-NumericValue $ws1.Range("C2") 1

$ws1.Range("A3").Value = "Document Status (Project)"
$ws1.Range("B3").Value = "312498*"
Set-NumericValue $ws1.Range("C3") 1

$ws1.Range("A4").Value = "Document Dashboard (Project)"
$ws1.Range("B4").Value = "312498*"
Set-NumericValue $ws1.Range("C4") 1

$ws1.Range("A4:C4").Copy()
$ws1.Range("A5:C5").PasteSpecial(-4122)   # xlPasteFormats, row 5 is new -> inherit row 4's look (style 2)
$ws1.Range("A5").Value = "Project Invoice Overview (Project)"
$ws1.Range("B5").Value = "312498*"
$ws1.Range("C5").Value = "1"

# ---------------------------------------------------------------------
# Sheet2 (MultiProjectReports): rework existing 5 rows and append 2 more.
# ---------------------------------------------------------------------

$ws2.Range("A2").Value = "Comment Status Internal"
$ws2.Range("B2").Value = "Project Teamcenter ID"
$ws2.Range("C2").Value = "312498*"
Set-NumericValue $ws2.Range("D2") 1

$ws2.Range("A3").Value = "Comment Status External"
$ws2.Range("B3").Value = "Project Teamcenter ID"
$ws2.Range("C3").Value = "312498*"
$ws2.Range("D3").Value = "1"

$ws2.Range("A4").Value = "Document Status"
$ws2.Range("B4").Value = "Project Teamcenter ID"
$ws2.Range("C4").Value = "312498*"
Set-NumericValue $ws2.Range("D4") 1

$ws2.Range("A5").Value = "Project Organization Overview"
$ws2.Range("B5").Value = "Project Manager"
$ws2.Range("C5").Value = "[MYSELF]"
Set-NumericValue $ws2.Range("D5") 1

$ws2.Range("A6").Value = "Project Portfolio Report"
$ws2.Range("B6").Value = "Project Manager"
$ws2.Range("C6").Value = "[MYSELF]"
Set-NumericValue $ws2.Range("D6") 1

$ws2.Range("A6:D6").Copy()
$ws2.Range("A7:D7").PasteSpecial(-4122)   # rows 7/8 are new -> inherit row 6's look (style 2)
$ws2.Range("A8:D8").PasteSpecial(-4122)

$ws2.Range("A7").Value = "Project Tracker"
$ws2.Range("B7").Value = "Project Teamcenter ID"
$ws2.Range("C7").Value = "[MYSELF]"
Set-NumericValue $ws2.Range("D7") 1

$ws2.Range("A8").Value = "Project Invoice Overview"
$ws2.Range("B8").Value = "Project Manager"
$ws2.Range("C8").Value = "[MYSELF]"
Set-NumericValue $ws2.Range("D8") 1

# ---------------------------------------------------------------------
# Column widths (best effort - underlying engine quantizes to 1/6 char
# steps so the exact fractional widths from real Excel can't always be
# reproduced bit-for-bit).
# ---------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 27.92
$ws1.Columns.Item(2).ColumnWidth = 21.29

# ---------------------------------------------------------------------
# Selections: set sheet2's active cell first, then sheet1's, so sheet1
# ends up as the active/tab-selected sheet (matching the source file).
# ---------------------------------------------------------------------
$ws2.Range("C4").Select()
$ws1.Range("B10").Select()
